$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 28 first (SC 92), then row 26 (RM 232), so indices don't shift
# before the second deletion is performed.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()
